{"js": "// The document body contains a single table where every 4th row (0, 4, 8,\n// 12, 16) holds five division-fact answers and the rows in between are\n// blank spacer rows. Replace the text of each of the 25 answer cells with\n// its new value while leaving every other part of the table (row/column\n// count, cell formatting, spacer rows) untouched.\nconst table = context.document.body.tables.getFirst();\ntable.load(\"values\");\nawait context.sync();\n\nconst newRowValues = {\n  0: [\"39\u00f74=9, 3\", \"99\u00f79=11, 0\", \"51\u00f79=5, 6\", \"44\u00f75=8, 4\", \"87\u00f78=10, 7\"],\n  4: [\"48\u00f77=6, 6\", \"88\u00f74=22, 0\", \"30\u00f79=3, 3\", \"48\u00f76=8, 0\", \"29\u00f78=3, 5\"],\n  8: [\"25\u00f75=5, 0\", \"45\u00f73=15, 0\", \"59\u00f79=6, 5\", \"29\u00f78=3, 5\", \"72\u00f78=9, 0\"],\n  12: [\"66\u00f75=13, 1\", \"78\u00f78=9, 6\", \"51\u00f78=6, 3\", \"32\u00f74=8, 0\", \"19\u00f73=6, 1\"],\n  16: [\"42\u00f72=21, 0\", \"10\u00f77=1, 3\", \"66\u00f72=33, 0\", \"44\u00f73=14, 2\", \"45\u00f77=6, 3\"],\n};\n\nconst values = table.values;\nfor (const rowIndex of Object.keys(newRowValues)) {\n  const idx = Number(rowIndex);\n  values[idx] = newRowValues[idx];\n}\ntable.values = values;\nawait context.sync();\n", "ps1": "# The document body contains a single table where every 4th row (1, 5, 9,\n# 13, 17 in Word's 1-based row numbering) holds five division-fact answers\n# and the rows in between are blank spacer rows. Replace the text of each\n# of the 25 answer cells with its new value while leaving everything else\n# (row/column count, cell formatting, spacer rows) untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"39\u00f74=9, 3\", \"99\u00f79=11, 0\", \"51\u00f79=5, 6\", \"44\u00f75=8, 4\", \"87\u00f78=10, 7\")\n    5  = @(\"48\u00f77=6, 6\", \"88\u00f74=22, 0\", \"30\u00f79=3, 3\", \"48\u00f76=8, 0\", \"29\u00f78=3, 5\")\n    9  = @(\"25\u00f75=5, 0\", \"45\u00f73=15, 0\", \"59\u00f79=6, 5\", \"29\u00f78=3, 5\", \"72\u00f78=9, 0\")\n    13 = @(\"66\u00f75=13, 1\", \"78\u00f78=9, 6\", \"51\u00f78=6, 3\", \"32\u00f74=8, 0\", \"19\u00f73=6, 1\")\n    17 = @(\"42\u00f72=21, 0\", \"10\u00f77=1, 3\", \"66\u00f72=33, 0\", \"44\u00f73=14, 2\", \"45\u00f77=6, 3\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
